$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B59").Value = 6662
$ws.Range("C59").Value = 5290
$ws.Range("D59").Value = 16504
$ws.Range("E59").Value = 11214
$ws.Range("F59").Value = 1372
$ws.Range("G59").Value = 2332
$ws.Range("H59").Value = 960
$ws.Range("I59").Value = 31010
$ws.Range("J59").Value = 27719
$ws.Range("K59").Value = 130
$ws.Range("L59").Value = 1427
$ws.Range("M59").Value = 4831
$ws.Range("N59").Value = 2373
$ws.Range("O59").Value = 2711
$ws.Range("P59").Value = 1618
$ws.Range("Q59").Value = 14628
$ws.Range("R59").Value = 3292
